$wb = $excel.ActiveWorkbook

# ===== Sheet: Mean Squared Error =====
$ws = $wb.Worksheets.Item("Mean Squared Error")
$ws.Range("H2").Value = 0.1139317528165474
$ws.Range("I2").Value = 0.1323295426629932
$ws.Range("J2").Value = 0.1499336727173778
$ws.Range("K2").Value = 0.168830334433119
$ws.Range("L2").Value = 0.1786301431989025
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("B4").Value = 0.003870036319383817
$ws.Range("C4").Value = 0.03722642358001035
$ws.Range("D4").Value = 0.07481188373261979
$ws.Range("E4").Value = 0.1125938691670729
$ws.Range("F4").Value = 0.1502775919606043
$ws.Range("G4").Value = 0.1887322985521378
$ws.Range("H4").Value = 0.2241673186014959
$ws.Range("I4").Value = 0.2608230629976908
$ws.Range("J4").Value = 0.3005325079271496
$ws.Range("K4").Value = 0.3344633537039347
$ws.Range("L4").Value = 0.373167496230116
$ws.Range("D5").Value = 0.03781977875788464
$ws.Range("E5").Value = 0.05630071391587962
$ws.Range("F5").Value = 0.07593812477087798
$ws.Range("J5").Value = 0.15009240467579
$ws.Range("K5").Value = 0.1686375884836185
$ws.Range("L5").Value = 0.1792612917002082
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("C7").Value = 0.01913098032101649
$ws.Range("I7").Value = 0.1322917493395617
$ws.Range("K7").Value = 0.1654364939889719
$ws.Range("L7").Value = 0.1675340234394192
$ws.Range("C8").Value = 0.01913098032101649
$ws.Range("J8").Value = 0.1463584243207595
$ws.Range("K8").Value = 0.1651870580543241
$ws.Range("L8").Value = 0.1708825118954486
$ws.Range("K9").Value = 0.1660033938404442
$ws.Range("L9").Value = 0.1677985767034396
$ws.Range("F10").Value = 0.07594946276790743
$ws.Range("G10").Value = 16231.70793319652
$ws.Range("H10").Value = 0.1139430908135769
$ws.Range("I10").Value = 0.1323182046659637
$ws.Range("J10").Value = 0.1501037426728194
$ws.Range("K10").Value = 16231.78326662812
$ws.Range("L10").Value = 16231.7956741762
$ws.Range("L11").Value = 0.1802514767741131
$ws.Range("I12").Value = 0.1323295426629932
$ws.Range("J12").Value = 0.1500810666787605
$ws.Range("K12").Value = 0.1540418069743799
$ws.Range("L12").Value = 0.1384482817265502
$ws.Range("K13").Value = 0.1697109188690726
$ws.Range("L13").Value = 0.1824964001859432
$ws.Range("I14").Value = 0.1320687687313159
$ws.Range("J14").Value = 0.1420348681201979
$ws.Range("K14").Value = 0.1434521177488785
$ws.Range("L14").Value = 0.1558407691697185
$ws.Range("L15").Value = 0.1855425420545206
$ws.Range("K16").Value = 0.1705499306492515
$ws.Range("L16").Value = 0.1833769846218967

# ===== Sheet: SNR =====
$ws = $wb.Worksheets.Item("SNR")
$ws.Range("H2").Value = 99.85647840542357
$ws.Range("I2").Value = 99.20635807945658
$ws.Range("J2").Value = 98.66393440071346
$ws.Range("K2").Value = 98.14842138747352
$ws.Range("L2").Value = 97.90337873124719
$ws.Range("B3").Value = "infinite"
$ws.Range("C3").Value = "infinite"
$ws.Range("D3").Value = "infinite"
$ws.Range("E3").Value = "infinite"
$ws.Range("F3").Value = "infinite"
$ws.Range("G3").Value = "infinite"
$ws.Range("H3").Value = "infinite"
$ws.Range("I3").Value = "infinite"
$ws.Range("J3").Value = "infinite"
$ws.Range("K3").Value = "infinite"
$ws.Range("L3").Value = "infinite"
$ws.Range("B4").Value = 114.7379172298362
$ws.Range("C4").Value = 104.9065544912582
$ws.Range("D4").Value = 101.8753617353521
$ws.Range("E4").Value = 100.0999202032218
$ws.Range("F4").Value = 98.84612536436006
$ws.Range("G4").Value = 97.85660534606841
$ws.Range("H4").Value = 97.10934466779599
$ws.Range("I4").Value = 96.4516077299519
$ws.Range("J4").Value = 95.83615308234175
$ws.Range("K4").Value = 95.3715822360107
$ws.Range("L4").Value = 94.89602955126401
$ws.Range("D5").Value = 104.5805115866948
$ws.Range("E5").Value = 102.8525624068806
$ws.Range("F5").Value = 101.5531027421092
$ws.Range("J5").Value = 98.5941142698244
$ws.Range("K5").Value = 98.08815759435753
$ws.Range("L5").Value = 97.82283621117379
$ws.Range("B6").Value = "infinite"
$ws.Range("C6").Value = "infinite"
$ws.Range("D6").Value = "infinite"
$ws.Range("E6").Value = "infinite"
$ws.Range("F6").Value = "infinite"
$ws.Range("G6").Value = "infinite"
$ws.Range("H6").Value = "infinite"
$ws.Range("I6").Value = "infinite"
$ws.Range("J6").Value = "infinite"
$ws.Range("K6").Value = "infinite"
$ws.Range("L6").Value = "infinite"
$ws.Range("C7").Value = 107.5383103876132
$ws.Range("I7").Value = 99.14035504414201
$ws.Range("K7").Value = 98.16936946003231
$ws.Range("L7").Value = 98.11465245247463
$ws.Range("C8").Value = 107.5342798180136
$ws.Range("J8").Value = 98.69747481502239
$ws.Range("K8").Value = 98.17189188218958
$ws.Range("L8").Value = 98.02467587446381
$ws.Range("K9").Value = 98.08116405145434
$ws.Range("L9").Value = 98.03445099420715
$ws.Range("F10").Value = 101.7330279032736
$ws.Range("G10").Value = 48.43463276958698
$ws.Range("H10").Value = 99.97139500607605
$ws.Range("I10").Value = 99.32207896927066
$ws.Range("J10").Value = 98.77435975497946
$ws.Range("K10").Value = 48.43461261347164
$ws.Range("L10").Value = 48.43460929373313
$ws.Range("L11").Value = 97.78237165581149
$ws.Range("I12").Value = 99.09918998570696
$ws.Range("J12").Value = 98.55249902400476
$ws.Range("K12").Value = 98.43937205694347
$ws.Range("L12").Value = 98.90288240017716
$ws.Range("K13").Value = 98.07758921197852
$ws.Range("L13").Value = 97.76214403827561
$ws.Range("I14").Value = 99.10264991110512
$ws.Range("J14").Value = 98.78670147706367
$ws.Range("K14").Value = 98.74358156056888
$ws.Range("L14").Value = 98.38384037021207
$ws.Range("L15").Value = 97.78746835128679
$ws.Range("K16").Value = 98.06548583061891
$ws.Range("L16").Value = 97.75055303121319

# ===== Sheet: PSNR =====
$ws = $wb.Worksheets.Item("PSNR")
$ws.Range("H2").Value = 105.7630182856831
$ws.Range("I2").Value = 105.1128979597161
$ws.Range("J2").Value = 104.570474280973
$ws.Range("K2").Value = 104.054961267733
$ws.Range("L2").Value = 103.8099186115067
$ws.Range("B3").Value = "infinite"
$ws.Range("C3").Value = "infinite"
$ws.Range("D3").Value = "infinite"
$ws.Range("E3").Value = "infinite"
$ws.Range("F3").Value = "infinite"
$ws.Range("G3").Value = "infinite"
$ws.Range("H3").Value = "infinite"
$ws.Range("I3").Value = "infinite"
$ws.Range("J3").Value = "infinite"
$ws.Range("K3").Value = "infinite"
$ws.Range("L3").Value = "infinite"
$ws.Range("B4").Value = 120.4523156674054
$ws.Range("C4").Value = 110.6209529288274
$ws.Range("D4").Value = 107.5897601729213
$ws.Range("E4").Value = 105.814318640791
$ws.Range("F4").Value = 104.5605238019292
$ws.Range("G4").Value = 103.5710037836376
$ws.Range("H4").Value = 102.8237431053652
$ws.Range("I4").Value = 102.166006167521
$ws.Range("J4").Value = 101.5505515199109
$ws.Range("K4").Value = 101.0859806735799
$ws.Range("L4").Value = 100.6104279888332
$ws.Range("D5").Value = 110.5522762359554
$ws.Range("E5").Value = 108.8243270561412
$ws.Range("F5").Value = 107.5248673913698
$ws.Range("J5").Value = 104.565878919085
$ws.Range("K5").Value = 104.0599222436181
$ws.Range("L5").Value = 103.7946008604344
$ws.Range("B6").Value = "infinite"
$ws.Range("C6").Value = "infinite"
$ws.Range("D6").Value = "infinite"
$ws.Range("E6").Value = "infinite"
$ws.Range("F6").Value = "infinite"
$ws.Range("G6").Value = "infinite"
$ws.Range("H6").Value = "infinite"
$ws.Range("I6").Value = "infinite"
$ws.Range("J6").Value = "infinite"
$ws.Range("K6").Value = "infinite"
$ws.Range("L6").Value = "infinite"
$ws.Range("C7").Value = 113.5120938255811
$ws.Range("I7").Value = 105.1141384821099
$ws.Range("K7").Value = 104.1431528980002
$ws.Range("L7").Value = 104.0884358904425
$ws.Range("C8").Value = 113.5120938255811
$ws.Range("J8").Value = 104.6752888225899
$ws.Range("K8").Value = 104.1497058897571
$ws.Range("L8").Value = 104.0024898820314
$ws.Range("K9").Value = 104.1282964050785
$ws.Range("L9").Value = 104.0815833478313
$ws.Range("F10").Value = 107.5242190132325
$ws.Range("G10").Value = 54.22582387954586
$ws.Range("H10").Value = 105.7625861160349
$ws.Range("I10").Value = 105.1132700792296
$ws.Range("J10").Value = 104.5655508649383
$ws.Range("K10").Value = 54.22580372343052
$ws.Range("L10").Value = 54.22580040369201
$ws.Range("L11").Value = 103.7706777601298
$ws.Range("I12").Value = 105.1128979597161
$ws.Range("J12").Value = 104.5662069980139
$ws.Range("K12").Value = 104.4530800309526
$ws.Range("L12").Value = 104.9165903741863
$ws.Range("K13").Value = 104.0323682265515
$ws.Range("L13").Value = 103.7169230528485
$ws.Range("I14").Value = 105.1214647857051
$ws.Range("J14").Value = 104.8055163516636
$ws.Range("K14").Value = 104.7623964351688
$ws.Range("L14").Value = 104.402655244812
$ws.Range("L15").Value = 103.6450310510904
$ws.Range("K16").Value = 104.0109506038429
$ws.Range("L16").Value = 103.6960178044372
